$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new handoff-file row ("cb15e4ab-bbac-446c-98e3-d44441351bb2.md") to
# each of the three report sheets: Overview, zh-cn, de-de. Each sheet has a
# ListObject (table) whose range needs to grow by one row.
# ---------------------------------------------------------------------------

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79cafe29fd773a2b1cc4303056992a99995f6d2a/e2e/"

# ---- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$rOverview = $loOverview.ListRows.Count

$wsOverview.Cells.Item($rOverview + 1, 1).Value = "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"
$wsOverview.Cells.Item($rOverview + 1, 2).Value = "e2e\cb15e4ab-bbac-446c-98e3-d44441351bb2.md"
$wsOverview.Cells.Item($rOverview + 1, 2).Style = $wsOverview.Cells.Item($rOverview, 2).Style
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rOverview + 1, 2), ($baseUrl + "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"), "", "", "e2e\cb15e4ab-bbac-446c-98e3-d44441351bb2.md") | Out-Null
$wsOverview.Cells.Item($rOverview + 1, 3).Value = ".md"
$wsOverview.Cells.Item($rOverview + 1, 4).Value = ""
$wsOverview.Cells.Item($rOverview + 1, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview + 1, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview + 1, 7).Value = "2016-09-05 08:54:24"
$wsOverview.Cells.Item($rOverview + 1, 7).Style = $wsOverview.Cells.Item($rOverview, 7).Style

# ---- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$rZh = $loZh.ListRows.Count

$wsZh.Cells.Item($rZh + 1, 1).Value = "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"
$wsZh.Cells.Item($rZh + 1, 1).Style = $wsZh.Cells.Item($rZh, 1).Style
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh + 1, 1), ($baseUrl + "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"), "", "", "cb15e4ab-bbac-446c-98e3-d44441351bb2.md") | Out-Null
$wsZh.Cells.Item($rZh + 1, 2).Value = ".md"
$wsZh.Cells.Item($rZh + 1, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($rZh + 1, 4).Value = "e2e"
$wsZh.Cells.Item($rZh + 1, 5).Value = "ht"
$wsZh.Cells.Item($rZh + 1, 6).Value = "False"
$wsZh.Cells.Item($rZh + 1, 7).Value = "cb15e4ab-bbac-446c-98e3-d44441351bb2.c4424cadd526157d1ba881655442df8db6fb82c7.zh-cn.xlf"
$wsZh.Cells.Item($rZh + 1, 8).Value = "2016-09-05 08:54:19"
$wsZh.Cells.Item($rZh + 1, 8).Style = $wsZh.Cells.Item($rZh, 8).Style
$wsZh.Cells.Item($rZh + 1, 9).Value = ""
$wsZh.Cells.Item($rZh + 1, 10).Value = ""
$wsZh.Cells.Item($rZh + 1, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($rZh + 1, 11).Style = $wsZh.Cells.Item($rZh, 11).Style
$wsZh.Cells.Item($rZh + 1, 12).Value = ""
$wsZh.Cells.Item($rZh + 1, 13).Value = "True"
$wsZh.Cells.Item($rZh + 1, 14).Value = ""
$wsZh.Cells.Item($rZh + 1, 15).Value = "False"
$wsZh.Cells.Item($rZh + 1, 16).Value = ""

# ---- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$rDe = $loDe.ListRows.Count

$wsDe.Cells.Item($rDe + 1, 1).Value = "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"
$wsDe.Cells.Item($rDe + 1, 1).Style = $wsDe.Cells.Item($rDe, 1).Style
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe + 1, 1), ($baseUrl + "cb15e4ab-bbac-446c-98e3-d44441351bb2.md"), "", "", "cb15e4ab-bbac-446c-98e3-d44441351bb2.md") | Out-Null
$wsDe.Cells.Item($rDe + 1, 2).Value = ".md"
$wsDe.Cells.Item($rDe + 1, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($rDe + 1, 4).Value = "e2e"
$wsDe.Cells.Item($rDe + 1, 5).Value = "ht"
$wsDe.Cells.Item($rDe + 1, 6).Value = "False"
$wsDe.Cells.Item($rDe + 1, 7).Value = "cb15e4ab-bbac-446c-98e3-d44441351bb2.c4424cadd526157d1ba881655442df8db6fb82c7.de-de.xlf"
$wsDe.Cells.Item($rDe + 1, 8).Value = "2016-09-05 08:54:24"
$wsDe.Cells.Item($rDe + 1, 8).Style = $wsDe.Cells.Item($rDe, 8).Style
$wsDe.Cells.Item($rDe + 1, 9).Value = ""
$wsDe.Cells.Item($rDe + 1, 10).Value = ""
$wsDe.Cells.Item($rDe + 1, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($rDe + 1, 11).Style = $wsDe.Cells.Item($rDe, 11).Style
$wsDe.Cells.Item($rDe + 1, 12).Value = ""
$wsDe.Cells.Item($rDe + 1, 13).Value = "True"
$wsDe.Cells.Item($rDe + 1, 14).Value = ""
$wsDe.Cells.Item($rDe + 1, 15).Value = "False"
$wsDe.Cells.Item($rDe + 1, 16).Value = ""

Write-Host "Done adding handback row."
